$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.011.12'
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").Value = '3.605.60'
$ws.Range("E3").Value = '  +1.17%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.78%  '
$ws.Range("D7").Value = '3.593.14'
$ws.Range("E7").Value = '  +1.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.606'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.195'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.42%  '
$ws.Range("E11").Value = '  +20.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.601'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '48.17'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000282'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("D15").Value = '4.183.53'
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '666.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '3.601.95'
$ws.Range("E18").Value = '  +2.36%  '
$ws.Range("D19").Value = '69.997.36'
$ws.Range("E19").Value = '  -1.92%  '
$ws.Range("E20").Value = '  -0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.925'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.90'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.91'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '576.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.99'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.106'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").Value = '3.561.73'
$ws.Range("E41").Value = '  -2.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0452'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.140'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.341'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.87%  '
$ws.Range("D46").Value = '0.0₃0731'
$ws.Range("E46").Value = '  -3.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.99%  '
$ws.Range("E48").Value = '  +4.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.132'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '135.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.33%  '
